# Weekly update: insert 4 new price records at the top of the
# "Vega Monumental Concepción - Zapallo" block (rows 121-124), pushing the
# existing rows down by 4 (121-138 -> 125-142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before the current row 121 so the existing data
# (formerly rows 121-138) shifts down to rows 125-142.
$ws.Rows("121:124").Insert()

# Row 121: Camote, 1a nueva(o)
$ws.Range("A121").Value = 11
$ws.Range("B121").Value = 'Vega Monumental Concepción'
$ws.Range("C121").Value = 'Bíobío'
$ws.Range("D121").Value = 44524
$ws.Range("E121").Value = 8
$ws.Range("F121").Value = 100112045
$ws.Range("G121").Value = 'Zapallo'
$ws.Range("H121").Value = 'Camote'
$ws.Range("I121").Value = '1a nueva(o)'
$ws.Range("J121").Value = 600
$ws.Range("K121").Value = 700
$ws.Range("L121").Value = 750
$ws.Range("M121").Value = 725
$ws.Range("N121").Value = '$/kilo (volumen en unidades)'
$ws.Range("O121").Value = 'Perú'
$ws.Range("P121").Value = 725
$ws.Range("Q121").Value = 1
$ws.Range("R121").Value = 'Hortaliza'

# Row 122: Camote, 2a nueva(o)
$ws.Range("A122").Value = 11
$ws.Range("B122").Value = 'Vega Monumental Concepción'
$ws.Range("C122").Value = 'Bíobío'
$ws.Range("D122").Value = 44524
$ws.Range("E122").Value = 8
$ws.Range("F122").Value = 100112045
$ws.Range("G122").Value = 'Zapallo'
$ws.Range("H122").Value = 'Camote'
$ws.Range("I122").Value = '2a nueva(o)'
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 650
$ws.Range("L122").Value = 650
$ws.Range("M122").Value = 650
$ws.Range("N122").Value = '$/kilo (volumen en unidades)'
$ws.Range("O122").Value = 'Perú'
$ws.Range("P122").Value = 650
$ws.Range("Q122").Value = 1
$ws.Range("R122").Value = 'Hortaliza'

# Row 123: Paine, 1a nueva(o)
$ws.Range("A123").Value = 11
$ws.Range("B123").Value = 'Vega Monumental Concepción'
$ws.Range("C123").Value = 'Bíobío'
$ws.Range("D123").Value = 44524
$ws.Range("E123").Value = 8
$ws.Range("F123").Value = 100112045
$ws.Range("G123").Value = 'Zapallo'
$ws.Range("H123").Value = 'Paine'
$ws.Range("I123").Value = '1a nueva(o)'
$ws.Range("J123").Value = 400
$ws.Range("K123").Value = 150
$ws.Range("L123").Value = 160
$ws.Range("M123").Value = 155
$ws.Range("N123").Value = '$/kilo (volumen en unidades)'
$ws.Range("O123").Value = 'Región de O''Higgins'
$ws.Range("P123").Value = 155
$ws.Range("Q123").Value = 1
$ws.Range("R123").Value = 'Hortaliza'

# Row 124: Paine, 2a nueva(o)
$ws.Range("A124").Value = 11
$ws.Range("B124").Value = 'Vega Monumental Concepción'
$ws.Range("C124").Value = 'Bíobío'
$ws.Range("D124").Value = 44524
$ws.Range("E124").Value = 8
$ws.Range("F124").Value = 100112045
$ws.Range("G124").Value = 'Zapallo'
$ws.Range("H124").Value = 'Paine'
$ws.Range("I124").Value = '2a nueva(o)'
$ws.Range("J124").Value = 200
$ws.Range("K124").Value = 100
$ws.Range("L124").Value = 100
$ws.Range("M124").Value = 100
$ws.Range("N124").Value = '$/kilo (volumen en unidades)'
$ws.Range("O124").Value = 'Región de O''Higgins'
$ws.Range("P124").Value = 100
$ws.Range("Q124").Value = 1
$ws.Range("R124").Value = 'Hortaliza'
